{"js": "// Locate the \"Registrar Viaje: ...\" bullet (last item of the \"Decisiones\n// tomadas en el dise\u00f1o de la aplicaci\u00f3n\" list) and add two new bullets\n// right after it, describing the \"Rendiciones\" (expense-report) query\n// assumptions, matching the existing list's style/numbering.\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  const t = paras.items[i].text;\n  if (t.indexOf(\"Registrar Viaje:\") !== -1 && t.indexOf(\"autom\u00f3vil\") !== -1) {\n    target = paras.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"No se encontr\u00f3 el p\u00e1rrafo 'Registrar Viaje' a continuar.\");\n}\n\nconst p1 = target.insertParagraph(\n  \"Rendiciones: tomamos en cuenta la hora fin del viaje ya que consideramos que se cobra una vez que se finalizo el viaje sin importar la fecha de inicio\",\n  Word.InsertLocation.after\n);\n\nconst p2 = p1.insertParagraph(\n  \"Rendiciones: el numero de rendiciones suponemos que es uno m\u00e1s que el de la ultima rendici\u00f3n y se genera autom\u00e1ticamente, es decir, el usuario no lo ingresa\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Locate the \"Registrar Viaje: ...\" bullet (last item of the \"Decisiones\n# tomadas en el dise\u00f1o de la aplicaci\u00f3n\" list) and add two new bullets\n# right after it, describing the \"Rendiciones\" (expense-report) query\n# assumptions, matching the existing list's style/numbering.\n$d = $word.ActiveDocument\n\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -like \"*Registrar Viaje:*\" -and $t -like \"*autom\u00f3vil*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"No se encontr\u00f3 el p\u00e1rrafo 'Registrar Viaje' a continuar.\"\n}\n\n$target = $d.Paragraphs.Item($targetIndex)\n$target.Range.InsertParagraphAfter()\n\n$p1 = $d.Paragraphs.Item($targetIndex + 1)\n$p1.Range.Text = \"Rendiciones: tomamos en cuenta la hora fin del viaje ya que consideramos que se cobra una vez que se finalizo el viaje sin importar la fecha de inicio\"\n\n$p1.Range.InsertParagraphAfter()\n\n$p2 = $d.Paragraphs.Item($targetIndex + 2)\n$p2.Range.Text = \"Rendiciones: el numero de rendiciones suponemos que es uno m\u00e1s que el de la ultima rendici\u00f3n y se genera autom\u00e1ticamente, es decir, el usuario no lo ingresa\"\n"}
